$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new blank column before column N (14th column), shifting
# the existing "Late" / "Date" / "Outstanding" columns right by one.
$ws.Columns("N:N").Insert()

# The newly inserted column inherits the width of the column
# immediately to its left (column M).
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Update the selection on the now-active "Repayment schedule" sheet.
$ws.Range("L20").Select() | Out-Null
